$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "[row:list ...]" loop marker (A6) now exposes a loop index variable,
# and the per-row id cell (A7) is switched from the data id to that index.
# NOTE: use single-quoted strings so PowerShell does not try to expand
# "${...}" as a variable reference.
$ws.Range("A6").Value = '[row:list datalist as data, index]'
$ws.Range("A7").Value = '${index}'

# Move the active selection to the (now) first data row of the list, A7:B7.
$ws.Range("A7:B7").Select()
